$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.913.32"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").Value = "3.537.39"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "608.00"
$ws.Range("E5").Value = "  +3.05%  "

$ws.Range("D6").Value = "185.51"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("D7").Value = "3.531.80"
$ws.Range("E7").Value = "  -0.97%  "

$ws.Range("D8").Value = "0.614"
$ws.Range("E8").Value = "  -1.37%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "0.213"
$ws.Range("E10").Value = "  +6.52%  "

$ws.Range("D11").Value = "0.641"
$ws.Range("E11").Value = "  -1.61%  "

$ws.Range("D12").Value = "53.65"
$ws.Range("E12").Value = "  -2.31%  "

$ws.Range("D13").Value = "0.0000308"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").Value = "9.47"
$ws.Range("E14").Value = "  -1.64%  "

$ws.Range("D15").Value = "4.096.16"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").Value = "69.981.25"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").Value = "12.63"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").Value = "18.88"
$ws.Range("E18").Value = "  -3.13%  "

$ws.Range("D19").Value = "3.532.93"
$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("D20").Value = "573.67"
$ws.Range("E20").Value = "  +1.61%  "

$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").Value = "0.990"
$ws.Range("E22").Value = "  -3.24%  "

$ws.Range("D23").Value = "17.39"
$ws.Range("E23").Value = "  -2.16%  "

$ws.Range("D24").Value = "4.69"
$ws.Range("E24").Value = "  -0.43%  "

$ws.Range("D25").Value = "4.87"
$ws.Range("E25").Value = "  -2.11%  "

$ws.Range("D26").Value = "93.95"
$ws.Range("E26").Value = "  -2.09%  "

$ws.Range("E27").Value = "  -1.13%  "

$ws.Range("D28").Value = "10.98"
$ws.Range("E28").Value = "  -4.76%  "

$ws.Range("D29").Value = "9.39"
$ws.Range("E29").Value = "  +2.02%  "

$ws.Range("D30").Value = "32.07"
$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("D31").Value = "7.01"
$ws.Range("E31").Value = "  -4.74%  "

$ws.Range("D32").Value = "12.16"
$ws.Range("E32").Value = "  -3.21%  "

$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  -1.66%  "

$ws.Range("D34").Value = "63.14"
$ws.Range("E34").Value = "  -3.16%  "

$ws.Range("D35").Value = "3.27"
$ws.Range("E35").Value = "  -0.42%  "

$ws.Range("D36").Value = "3.60"
$ws.Range("E36").Value = "  +16.45%  "

$ws.Range("D37").Value = "535.63"
$ws.Range("E37").Value = "  -4.00%  "

$ws.Range("D38").Value = "0.402"
$ws.Range("E38").Value = "  -3.77%  "

$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.06%  "

$ws.Range("D40").Value = "37.09"
$ws.Range("E40").Value = "  -3.17%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0779"
$ws.Range("E41").Value = "  +0.23%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.530.04"
$ws.Range("E42").Value = "  +5.33%  "

$ws.Range("D43").Value = "3.52"
$ws.Range("E43").Value = "  +3.89%  "

$ws.Range("D44").Value = "0.136"
$ws.Range("E44").Value = "  +0.67%  "

$ws.Range("D45").Value = "0.0453"
$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("D46").Value = "2.93"
$ws.Range("E46").Value = "  -2.22%  "

$ws.Range("D47").Value = "3.41"
$ws.Range("E47").Value = "  -4.55%  "

$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  -4.21%  "

$ws.Range("E50").Value = "  +0.21%  "

$ws.Range("E51").Value = "  -4.04%  "

